$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated "close" values for existing rows 2-127 (index 0-125) ---
$closeValues = @(
  0.85915,
  0.85784,
  0.85425,
  0.91046,
  0.90915,
  0.8732,
  0.8924800000000001,
  0.88987,
  0.89118,
  0.89641,
  0.8951,
  0.8924800000000001,
  0.88235,
  0.89444,
  0.91895,
  0.90588,
  0.91732,
  0.91405,
  0.91895,
  0.91569,
  0.91569,
  0.91895,
  0.92222,
  0.91569,
  0.91569,
  0.92386,
  0.9385599999999999,
  0.96471,
  0.98431,
  0.98431,
  0.99085,
  1.06601,
  1.08562,
  1.08725,
  1.07745,
  1.04967,
  0.99902,
  0.99085,
  0.99575,
  0.97941,
  0.96797,
  0.95817,
  0.96797,
  0.95163,
  0.93529,
  0.95163,
  0.95163,
  0.9271200000000001,
  0.9271200000000001,
  0.92549,
  0.94837,
  0.94346,
  0.95163,
  0.94837,
  0.98105,
  1.01699,
  1.01209,
  0.98268,
  0.97614,
  0.97451,
  0.96634,
  0.96797,
  0.96634,
  0.95817,
  0.9549,
  0.94837,
  0.95163,
  0.95817,
  0.95163,
  0.95,
  0.99085,
  0.97778,
  0.99248,
  0.98431,
  0.96307,
  0.9451000000000001,
  0.92549,
  0.97124,
  1.01536,
  1.00556,
  1.00065,
  0.98268,
  0.9418299999999999,
  0.95327,
  0.9451000000000001,
  0.95163,
  0.92876,
  0.9565399999999999,
  0.9565399999999999,
  0.96634,
  0.95,
  0.95817,
  0.9565399999999999,
  1.04804,
  1.14608,
  1.1902,
  1.11667,
  1.21307,
  1.24248,
  1.21634,
  1.25556,
  1.14281,
  1.14281,
  1.08725,
  1.00229,
  1.02843,
  1.04314,
  1.12484,
  1.10196,
  1.10359,
  1.15261,
  1.1183,
  1.17549,
  1.22778,
  1.33725,
  1.4549,
  1.56928,
  1.55294,
  1.68366,
  1.70327,
  1.62484,
  1.64771,
  1.57255,
  1.65425,
  1.64444,
  1.72288
)

for ($i = 0; $i -lt $closeValues.Length; $i++) {
  $row = 2 + $i
  $ws.Cells.Item($row, 2).Value = $closeValues[$i]
}

# --- New rows 128-133 (index 126-131) appended to the series ---
$newRows = @(
  @(126, 1.59542),
  @(127, 1.57255),
  @(128, 1.67386),
  @(129, 1.60523),
  @(130, 1.57908),
  @(131, 1.63464)
)

$startRow = 128
for ($i = 0; $i -lt $newRows.Length; $i++) {
  $row = $startRow + $i
  $ws.Cells.Item($row, 1).Value = $newRows[$i][0]
  $ws.Cells.Item($row, 2).Value = $newRows[$i][1]
}

# Copy the index-column style (bold, centered, bordered) from row 127 down to the new rows
$ws.Range("A127").Copy()
$ws.Range("A128:A133").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Output "Updated close_validation_scaler data and appended 6 new rows"
